$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clean up / correct the weather-station names in column A (rows 11-36):
# fix typos, drop extraneous words ("airport", "forest", "park", "niab"),
# and strip spaces/hyphens so each name is a single token.
$ws.Range("A11").Value = "cardiff"
$ws.Range("A12").Value = "rossonwye"
$ws.Range("A13").Value = "aberporth"
$ws.Range("A14").Value = "cwmystwyth"
$ws.Range("A15").Value = "cambridge"
$ws.Range("A16").Value = "lowestoft"
$ws.Range("A17").Value = "shawbury"
$ws.Range("A18").Value = "suttonbonington"
$ws.Range("A19").Value = "waddington"
$ws.Range("A20").Value = "sheffield"
$ws.Range("A21").Value = "valley"
$ws.Range("A22").Value = "ringway"
$ws.Range("A23").Value = "bradford"
$ws.Range("A24").Value = "whitby"
$ws.Range("A25").Value = "newtonrigg"
$ws.Range("A26").Value = "durham"
$ws.Range("A27").Value = "eskdalemuir"
$ws.Range("A28").Value = "paisley"
$ws.Range("A29").Value = "leuchars"
$ws.Range("A30").Value = "dunstaffnage"
$ws.Range("A31").Value = "iree"
$ws.Range("A32").Value = "braemar"
$ws.Range("A33").Value = "stornoway"
$ws.Range("A34").Value = "wickairpor"
$ws.Range("A35").Value = "ballypatrick"
$ws.Range("A36").Value = "armagh"

# Widen column A so the (now different-length) station names are readable.
$ws.Columns.Item(1).ColumnWidth = 20.333333333333332

# Move the active selection to A34.
$ws.Range("A34").Select()
